# until error page in toseesaderat
# Add a second task row ("list of all bank branches") to the tracker sheet,
# wire up its formula/merges, and nudge column B wide enough for the new title.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 5: new task #2 ("بانک ملت" branch list task), mirrors row 2's layout ----
$ws.Range("A2").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("B2").Copy()
$ws.Range("B5").PasteSpecial(-4122)
$ws.Range("C2").Copy()
$ws.Range("C5").PasteSpecial(-4122)
$ws.Range("D2").Copy()
$ws.Range("D5").PasteSpecial(-4122)
$ws.Range("D2").Copy()
$ws.Range("E5").PasteSpecial(-4122)
$ws.Range("D2").Copy()
$ws.Range("F5").PasteSpecial(-4122)
$ws.Range("D2").Copy()
$ws.Range("G5").PasteSpecial(-4122)

$ws.Range("A5").Value = 2
$ws.Range("B5").Value = "لیست تمام شعب کل بانک های کشور"
$ws.Range("C5").Value = "تمام بانک های کشور مشخص شود و ابتدا یک لیست درست شود که از چه طریقی می توان به لیست شعب رسید مثلا خودشون فرم اکسل دارن یا نه یه جدوله که با پایتون میشه خوندش یا نه هچ شعبی نداره یا پیدا نمیشه و بعد ایچاد یک اکسل بزرگ اطلاعات به ترتیب باشه برای همه بانک ها"

# D5 repeats the same date text as D4 ("1403/08/07") - force text (not a Jalali date
# serial) the same way the sheet already stores its date column, then restore the
# plain (non-text-forced) number format that the rest of column D uses.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1403/08/07"
$ws.Range("D2").Copy()
$ws.Range("D5").PasteSpecial(-4122)

$ws.Range("E5").Value = 18.5
$ws.Range("F5").Value = 21
$ws.Range("G5").Formula = "=F5-E5"

$ws.Rows.Item(5).RowHeight = 57.6

# ---- Row 6: continuation line (second sub-task, "بانک ملت") ----
$ws.Range("A2").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("B2").Copy()
$ws.Range("B6").PasteSpecial(-4122)
$ws.Range("D2").Copy()
$ws.Range("C6").PasteSpecial(-4122)
$ws.Range("D2").Copy()
$ws.Range("D6").PasteSpecial(-4122)
$ws.Range("D2").Copy()
$ws.Range("E6").PasteSpecial(-4122)

$ws.Range("C6").Value = "بانک ملت"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1403/08/10"
$ws.Range("D2").Copy()
$ws.Range("D6").PasteSpecial(-4122)

$ws.Range("E6").Value = 1.5

# ---- merges for the new task block ----
$ws.Range("A5:A6").Merge()
$ws.Range("B5:B6").Merge()

# ---- column B needs to be a bit wider for the new (slightly longer) title ----
$ws.Columns.Item(2).ColumnWidth = 26.65

# ---- match the saved selection/cursor position ----
$ws.Range("F6").Select()

$wb.Application.CutCopyMode = $false
